$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.527.97"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "2.105.16"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "332.29"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.5232"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4500"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").Value = "53.63"
$ws.Range("E9").Value = "  +16.68%  "
$ws.Range("D10").Value = "0.08946"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").Value = "1.158"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "24.51"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "2.098.17"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "6.735"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "7.721"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "96.45"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001126"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").Value = "0.06628"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "19.22"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.300"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "30.554.88"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").Value = "12.36"
$ws.Range("D25").Value = "2.326"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").Value = "2.337.20"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "22.32"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "2.587"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "163.82"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.40"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "1.201"
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("D32").Value = "0.1075"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").Value = "1.675"
$ws.Range("E33").Value = "  +8.62%  "
$ws.Range("D34").Value = "6.165"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "3.904"
$ws.Range("E35").Value = "  -3.72%  "
$ws.Range("D36").Value = "10.48"
$ws.Range("E36").Value = "  +10.24%  "
$ws.Range("D37").Value = "0.02573"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "0.06785"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "5.494"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "12.76"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "0.2267"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "0.6929"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "1.257"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "0.6371"
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("D47").Value = "2.298"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("D48").Value = "3.641"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("D50").Value = "1.246"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("D51").Value = "82.15"
$ws.Range("E51").Value = "  -0.29%  "
